$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.180.90"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -5.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.557.65"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.93"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.04"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -5.63%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -2.87%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.05"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.72"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -4.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.110"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.946.93"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.542.36"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.876"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -3.68%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.173.34"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -6.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.27"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +5.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.30"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "261.14"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -9.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -2.87%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.79"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.33%  "

# Row 26
$ws.Range("E26").Value = "  -4.75%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -5.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.49"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -2.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.14"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.77"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.40"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -6.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0801"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -4.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.83%  "

# Row 38
$ws.Range("E38").Value = "  -2.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.87"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +8.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.51"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +12.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.43%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0314"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.91"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.081.84"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.02"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -9.73%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -4.31%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.58"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.94%  "

# Row 49
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.803.27"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.76"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.42%  "

# Row 51
$ws.Range("E51").Value = "  -1.68%  "
